# Updates the NBA fantasy roster table ("LOS Galacticos") on Sheet1.
# Column A (player) values/order are unchanged; Column B (position) and
# Column C (team) are refreshed to the new roster assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Fred VanVleet",            "PG",       "Houston Rockets"),
    @("Kentavious Caldwell-Pope", "SG,SF",    "Orlando Magic"),
    @("Anthony Edwards",          "SG,SF",    "Minnesota Timberwolves"),
    @("Anfernee Simons",          "PG,SG",    "Portland Trail Blazers"),
    @("Paul George",              "SG,SF,PF", "Philadelphia 76ers"),
    @("Kyle Kuzma",                "PF",       "Washington Wizards"),
    @("Jayson Tatum",              "SF,PF",    "Boston Celtics"),
    @("Jonathan Kuminga",          "SF,PF",    "Golden State Warriors"),
    @("Jaren Jackson Jr.",         "PF,C",     "Memphis Grizzlies"),
    @("Zach Edey",                 "C",        "Memphis Grizzlies"),
    @("Andrew Wiggins",            "SF,PF",    "Golden State Warriors"),
    @("Scoot Henderson",           "PG",       "Portland Trail Blazers"),
    @("Dennis Schröder",           "PG,SG",    "Golden State Warriors"),
    @("Ivica Zubac",               "C",        "LA Clippers"),
    @("Giannis Antetokounmpo",     "PF,C",     "Milwaukee Bucks"),
    @("Jaden Ivey",                "PG,SG",    "Detroit Pistons"),
    @("Zion Williamson",           "PF,C",     "New Orleans Pelicans"),
    @("James Harden",              "PG,SG",    "LA Clippers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
